$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.5984806666666667
$ws.Range("H2").Value = 1.795442
$ws.Range("I2").Value = 0.002763845209804117
$ws.Range("J2").Value = 0.002763845209804117
$ws.Range("M2").Value = 25.69910333333333
$ws.Range("N2").Value = 77.09731
$ws.Range("O2").Value = 0.08761243344445813
$ws.Range("P2").Value = 0.08761243344445814
$ws.Range("Q2").Value = 15.38041649566889
$ws.Range("R2").Value = 138.42374846102
$ws.Range("S2").Value = 0.0002421472044947476
$ws.Range("T2").Value = 0.0002421472044947477

$ws.Range("G3").Value = 0.5984806666666667
$ws.Range("H3").Value = 1.795442
$ws.Range("I3").Value = 0.002763845209804117
$ws.Range("J3").Value = 0.002763845209804117
$ws.Range("O3").Value = 0.0005530844306649811
$ws.Range("P3").Value = 0.0005530844306649812
$ws.Range("Q3").Value = 0.09709431146311111
$ws.Range("R3").Value = 0.8738488031680001
$ws.Range("S3").Value = 0.000001528639754310646
$ws.Range("T3").Value = 0.000001528639754310646

$ws.Range("G4").Value = 0.5984806666666667
$ws.Range("H4").Value = 1.795442
$ws.Range("I4").Value = 0.002763845209804117
$ws.Range("J4").Value = 0.002763845209804117
$ws.Range("M4").Value = 46.92720933333334
$ws.Range("N4").Value = 140.781628
$ws.Range("O4").Value = 0.1599825079935015
$ws.Range("P4").Value = 0.1599825079935016
$ws.Range("Q4").Value = 28.08502752661956
$ws.Range("R4").Value = 252.765247739576
$ws.Range("S4").Value = 0.0004421668883702881
$ws.Range("T4").Value = 0.0004421668883702882

$ws.Range("G5").Value = 0.5984806666666667
$ws.Range("H5").Value = 1.795442
$ws.Range("I5").Value = 0.002763845209804117
$ws.Range("J5").Value = 0.002763845209804117
$ws.Range("M5").Value = 220.538579
$ws.Range("N5").Value = 661.615737
$ws.Range("O5").Value = 0.7518519741313753
$ws.Range("P5").Value = 0.7518519741313754
$ws.Range("Q5").Value = 131.9880757856393
$ws.Range("R5").Value = 1187.892682070754
$ws.Range("S5").Value = 0.002078002477184771
$ws.Range("T5").Value = 0.002078002477184771

$ws.Range("I6").Value = 0.9924988388011551
$ws.Range("J6").Value = 0.9924988388011552
$ws.Range("M6").Value = 25.69910333333333
$ws.Range("N6").Value = 77.09731
$ws.Range("O6").Value = 0.08761243344445813
$ws.Range("P6").Value = 0.08761243344445814
$ws.Range("Q6").Value = 5523.118826655052
$ws.Range("R6").Value = 49708.06943989547
$ws.Range("S6").Value = 0.08695523845816817
$ws.Range("T6").Value = 0.0869552384581682

$ws.Range("I7").Value = 0.9924988388011551
$ws.Range("J7").Value = 0.9924988388011552
$ws.Range("O7").Value = 0.0005530844306649811
$ws.Range("P7").Value = 0.0005530844306649812
$ws.Range("S7").Value = 0.0005489356551939917
$ws.Range("T7").Value = 0.0005489356551939919

$ws.Range("I8").Value = 0.9924988388011551
$ws.Range("J8").Value = 0.9924988388011552
$ws.Range("M8").Value = 46.92720933333334
$ws.Range("N8").Value = 140.781628
$ws.Range("O8").Value = 0.1599825079935015
$ws.Range("P8").Value = 0.1599825079935016
$ws.Range("Q8").Value = 10085.35395144069
$ws.Range("R8").Value = 90768.18556296625
$ws.Range("S8").Value = 0.1587824534120468
$ws.Range("T8").Value = 0.1587824534120468

$ws.Range("I9").Value = 0.9924988388011551
$ws.Range("J9").Value = 0.9924988388011552
$ws.Range("M9").Value = 220.538579
$ws.Range("N9").Value = 661.615737
$ws.Range("O9").Value = 0.7518519741313753
$ws.Range("P9").Value = 0.7518519741313754
$ws.Range("Q9").Value = 47397.01466933098
$ws.Range("R9").Value = 426573.1320239787
$ws.Range("S9").Value = 0.7462122112757461
$ws.Range("T9").Value = 0.7462122112757462

$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.278701
$ws.Range("H10").Value = 0.836103
$ws.Range("I10").Value = 0.001287069853246639
$ws.Range("J10").Value = 0.001287069853246639
$ws.Range("M10").Value = 25.69910333333333
$ws.Range("N10").Value = 77.09731
$ws.Range("O10").Value = 0.08761243344445813
$ws.Range("P10").Value = 0.08761243344445814
$ws.Range("Q10").Value = 7.162365798103333
$ws.Range("R10").Value = 64.46129218293
$ws.Range("S10").Value = 0.0001127633218559397
$ws.Range("T10").Value = 0.0001127633218559397

$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 0.6666666666666666
$ws.Range("G11").Value = 0.278701
$ws.Range("H11").Value = 0.836103
$ws.Range("I11").Value = 0.001287069853246639
$ws.Range("J11").Value = 0.001287069853246639
$ws.Range("O11").Value = 0.0005530844306649811
$ws.Range("P11").Value = 0.0005530844306649812
$ws.Range("Q11").Value = 0.04521496383466667
$ws.Range("R11").Value = 0.406934674512
$ws.Range("S11").Value = 0.0000007118582970089782
$ws.Range("T11").Value = 0.0000007118582970089783

$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.278701
$ws.Range("H12").Value = 0.836103
$ws.Range("I12").Value = 0.001287069853246639
$ws.Range("J12").Value = 0.001287069853246639
$ws.Range("M12").Value = 46.92720933333334
$ws.Range("N12").Value = 140.781628
$ws.Range("O12").Value = 0.1599825079935015
$ws.Range("P12").Value = 0.1599825079935016
$ws.Range("Q12").Value = 13.07866016840934
$ws.Range("R12").Value = 117.707941515684
$ws.Range("S12").Value = 0.0002059086630852253
$ws.Range("T12").Value = 0.0002059086630852253

$ws.Range("E13").Value = 2
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.278701
$ws.Range("H13").Value = 0.836103
$ws.Range("I13").Value = 0.001287069853246639
$ws.Range("J13").Value = 0.001287069853246639
$ws.Range("M13").Value = 220.538579
$ws.Range("N13").Value = 661.615737
$ws.Range("O13").Value = 0.7518519741313753
$ws.Range("P13").Value = 0.7518519741313754
$ws.Range("Q13").Value = 61.464322505879
$ws.Range("R13").Value = 553.178902552911
$ws.Range("S13").Value = 0.000967686010008465
$ws.Range("T13").Value = 0.0009676860100084651

$ws.Range("G14").Value = 0.7471133333333334
$ws.Range("H14").Value = 2.24134
$ws.Range("I14").Value = 0.003450246135794061
$ws.Range("J14").Value = 0.003450246135794061
$ws.Range("M14").Value = 25.69910333333333
$ws.Range("N14").Value = 77.09731
$ws.Range("O14").Value = 0.08761243344445813
$ws.Range("P14").Value = 0.08761243344445814
$ws.Range("Q14").Value = 19.20014275504444
$ws.Range("R14").Value = 172.8012847954
$ws.Range("S14").Value = 0.000302284459939256
$ws.Range("T14").Value = 0.0003022844599392561

$ws.Range("G15").Value = 0.7471133333333334
$ws.Range("H15").Value = 2.24134
$ws.Range("I15").Value = 0.003450246135794061
$ws.Range("J15").Value = 0.003450246135794061
$ws.Range("O15").Value = 0.0005530844306649811
$ws.Range("P15").Value = 0.0005530844306649812
$ws.Range("Q15").Value = 0.1212076825955556
$ws.Range("R15").Value = 1.09086914336
$ws.Range("S15").Value = 0.00000190827741966971
$ws.Range("T15").Value = 0.00000190827741966971

$ws.Range("G16").Value = 0.7471133333333334
$ws.Range("H16").Value = 2.24134
$ws.Range("I16").Value = 0.003450246135794061
$ws.Range("J16").Value = 0.003450246135794061
$ws.Range("M16").Value = 46.92720933333334
$ws.Range("N16").Value = 140.781628
$ws.Range("O16").Value = 0.1599825079935015
$ws.Range("P16").Value = 0.1599825079935016
$ws.Range("Q16").Value = 35.05994378905778
$ws.Range("R16").Value = 315.5394941015201
$ws.Range("S16").Value = 0.0005519790299992212
$ws.Range("T16").Value = 0.0005519790299992213

$ws.Range("G17").Value = 0.7471133333333334
$ws.Range("H17").Value = 2.24134
$ws.Range("I17").Value = 0.003450246135794061
$ws.Range("J17").Value = 0.003450246135794061
$ws.Range("M17").Value = 220.538579
$ws.Range("N17").Value = 661.615737
$ws.Range("O17").Value = 0.7518519741313753
$ws.Range("P17").Value = 0.7518519741313754
$ws.Range("Q17").Value = 164.7673128852867
$ws.Range("R17").Value = 1482.90581596758
$ws.Range("S17").Value = 0.002594074368435914
$ws.Range("T17").Value = 0.002594074368435915
